$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.145036666666667
$ws.Cells.Item(2, 8).Value = 3.43511
$ws.Cells.Item(2, 9).Value = 0.4953865629219574
$ws.Cells.Item(2, 10).Value = 0.4953865629219574
$ws.Cells.Item(2, 13).Value = 41.83357366666667
$ws.Cells.Item(2, 14).Value = 125.500721
$ws.Cells.Item(2, 15).Value = 0.2773195847425811
$ws.Cells.Item(2, 16).Value = 0.2773195847425811
$ws.Cells.Item(2, 17).Value = 47.90097574603445
$ws.Cells.Item(2, 18).Value = 431.10878171431
$ws.Cells.Item(2, 19).Value = 0.1373803959165718
$ws.Cells.Item(2, 20).Value = 0.1373803959165718

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.145036666666667
$ws.Cells.Item(3, 8).Value = 3.43511
$ws.Cells.Item(3, 9).Value = 0.4953865629219574
$ws.Cells.Item(3, 10).Value = 0.4953865629219574
$ws.Cells.Item(3, 15).Value = 0.4239803668761465
$ws.Cells.Item(3, 16).Value = 0.4239803668761465
$ws.Cells.Item(3, 17).Value = 73.23346199793555
$ws.Cells.Item(3, 18).Value = 659.1011579814199
$ws.Cells.Item(3, 19).Value = 0.2100341766931647
$ws.Cells.Item(3, 20).Value = 0.2100341766931647

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.145036666666667
$ws.Cells.Item(4, 8).Value = 3.43511
$ws.Cells.Item(4, 9).Value = 0.4953865629219574
$ws.Cells.Item(4, 10).Value = 0.4953865629219574
$ws.Cells.Item(4, 13).Value = 14.18032733333333
$ws.Cells.Item(4, 14).Value = 42.540982
$ws.Cells.Item(4, 15).Value = 0.09400302539123752
$ws.Cells.Item(4, 16).Value = 0.09400302539123752
$ws.Cells.Item(4, 17).Value = 16.23699474200222
$ws.Cells.Item(4, 18).Value = 146.13295267802
$ws.Cells.Item(4, 19).Value = 0.04656783565283065
$ws.Cells.Item(4, 20).Value = 0.04656783565283065

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.145036666666667
$ws.Cells.Item(5, 8).Value = 3.43511
$ws.Cells.Item(5, 9).Value = 0.4953865629219574
$ws.Cells.Item(5, 10).Value = 0.4953865629219574
$ws.Cells.Item(5, 13).Value = 30.87848266666667
$ws.Cells.Item(5, 14).Value = 92.635448
$ws.Cells.Item(5, 15).Value = 0.204697022990035
$ws.Cells.Item(5, 16).Value = 0.2046970229900349
$ws.Cells.Item(5, 17).Value = 35.35699486436445
$ws.Cells.Item(5, 18).Value = 318.21295377928
$ws.Cells.Item(5, 19).Value = 0.1014041546593903
$ws.Cells.Item(5, 20).Value = 0.1014041546593903

# Row 6
$ws.Cells.Item(6, 7).Value = 0.4713496666666666
$ws.Cells.Item(6, 9).Value = 0.2039238551060172
$ws.Cells.Item(6, 10).Value = 0.2039238551060172
$ws.Cells.Item(6, 13).Value = 41.83357366666667
$ws.Cells.Item(6, 14).Value = 125.500721
$ws.Cells.Item(6, 15).Value = 0.2773195847425811
$ws.Cells.Item(6, 16).Value = 0.2773195847425811
$ws.Cells.Item(6, 17).Value = 19.71824100325878
$ws.Cells.Item(6, 18).Value = 177.464169029329
$ws.Cells.Item(6, 19).Value = 0.05655207881710698
$ws.Cells.Item(6, 20).Value = 0.05655207881710697

# Row 7
$ws.Cells.Item(7, 7).Value = 0.4713496666666666
$ws.Cells.Item(7, 9).Value = 0.2039238551060172
$ws.Cells.Item(7, 10).Value = 0.2039238551060172
$ws.Cells.Item(7, 15).Value = 0.4239803668761465
$ws.Cells.Item(7, 16).Value = 0.4239803668761465
$ws.Cells.Item(7, 17).Value = 30.14625549246421
$ws.Cells.Item(7, 19).Value = 0.08645971090264733
$ws.Cells.Item(7, 20).Value = 0.08645971090264733

# Row 8
$ws.Cells.Item(8, 7).Value = 0.4713496666666666
$ws.Cells.Item(8, 9).Value = 0.2039238551060172
$ws.Cells.Item(8, 10).Value = 0.2039238551060172
$ws.Cells.Item(8, 13).Value = 14.18032733333333
$ws.Cells.Item(8, 14).Value = 42.540982
$ws.Cells.Item(8, 15).Value = 0.09400302539123752
$ws.Cells.Item(8, 16).Value = 0.09400302539123752
$ws.Cells.Item(8, 17).Value = 6.683892561790888
$ws.Cells.Item(8, 18).Value = 60.15503305611799
$ws.Cells.Item(8, 19).Value = 0.01916945932940998
$ws.Cells.Item(8, 20).Value = 0.01916945932940998

# Row 9
$ws.Cells.Item(9, 7).Value = 0.4713496666666666
$ws.Cells.Item(9, 9).Value = 0.2039238551060172
$ws.Cells.Item(9, 10).Value = 0.2039238551060172
$ws.Cells.Item(9, 13).Value = 30.87848266666667
$ws.Cells.Item(9, 14).Value = 92.635448
$ws.Cells.Item(9, 15).Value = 0.204697022990035
$ws.Cells.Item(9, 16).Value = 0.2046970229900349
$ws.Cells.Item(9, 17).Value = 14.55456251210578
$ws.Cells.Item(9, 18).Value = 130.991062608952
$ws.Cells.Item(9, 19).Value = 0.04174260605685298
$ws.Cells.Item(9, 20).Value = 0.04174260605685297

# Row 10
$ws.Cells.Item(10, 7).Value = 0.6323219999999999
$ws.Cells.Item(10, 8).Value = 1.896966
$ws.Cells.Item(10, 9).Value = 0.2735666300991275
$ws.Cells.Item(10, 10).Value = 0.2735666300991275
$ws.Cells.Item(10, 13).Value = 41.83357366666667
$ws.Cells.Item(10, 14).Value = 125.500721
$ws.Cells.Item(10, 15).Value = 0.2773195847425811
$ws.Cells.Item(10, 16).Value = 0.2773195847425811
$ws.Cells.Item(10, 17).Value = 26.452288968054
$ws.Cells.Item(10, 18).Value = 238.070600712486
$ws.Cells.Item(10, 19).Value = 0.07586538425851733
$ws.Cells.Item(10, 20).Value = 0.07586538425851733

# Row 11
$ws.Cells.Item(11, 7).Value = 0.6323219999999999
$ws.Cells.Item(11, 8).Value = 1.896966
$ws.Cells.Item(11, 9).Value = 0.2735666300991275
$ws.Cells.Item(11, 10).Value = 0.2735666300991275
$ws.Cells.Item(11, 15).Value = 0.4239803668761465
$ws.Cells.Item(11, 16).Value = 0.4239803668761465
$ws.Cells.Item(11, 17).Value = 40.44161248762799
$ws.Cells.Item(11, 18).Value = 363.9745123886519
$ws.Cells.Item(11, 19).Value = 0.1159868801944991
$ws.Cells.Item(11, 20).Value = 0.1159868801944991

# Row 12
$ws.Cells.Item(12, 7).Value = 0.6323219999999999
$ws.Cells.Item(12, 8).Value = 1.896966
$ws.Cells.Item(12, 9).Value = 0.2735666300991275
$ws.Cells.Item(12, 10).Value = 0.2735666300991275
$ws.Cells.Item(12, 13).Value = 14.18032733333333
$ws.Cells.Item(12, 14).Value = 42.540982
$ws.Cells.Item(12, 15).Value = 0.09400302539123752
$ws.Cells.Item(12, 16).Value = 0.09400302539123752
$ws.Cells.Item(12, 17).Value = 8.966532940067999
$ws.Cells.Item(12, 18).Value = 80.698796460612
$ws.Cells.Item(12, 19).Value = 0.02571609087540356
$ws.Cells.Item(12, 20).Value = 0.02571609087540357

# Row 13
$ws.Cells.Item(13, 7).Value = 0.6323219999999999
$ws.Cells.Item(13, 8).Value = 1.896966
$ws.Cells.Item(13, 9).Value = 0.2735666300991275
$ws.Cells.Item(13, 10).Value = 0.2735666300991275
$ws.Cells.Item(13, 13).Value = 30.87848266666667
$ws.Cells.Item(13, 14).Value = 92.635448
$ws.Cells.Item(13, 15).Value = 0.204697022990035
$ws.Cells.Item(13, 16).Value = 0.2046970229900349
$ws.Cells.Item(13, 17).Value = 19.525143916752
$ws.Cells.Item(13, 18).Value = 175.726295250768
$ws.Cells.Item(13, 19).Value = 0.05599827477070749
$ws.Cells.Item(13, 20).Value = 0.05599827477070749

# Row 14
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.062692
$ws.Cells.Item(14, 8).Value = 0.188076
$ws.Cells.Item(14, 9).Value = 0.02712295187289783
$ws.Cells.Item(14, 10).Value = 0.02712295187289783
$ws.Cells.Item(14, 13).Value = 41.83357366666667
$ws.Cells.Item(14, 14).Value = 125.500721
$ws.Cells.Item(14, 15).Value = 0.2773195847425811
$ws.Cells.Item(14, 16).Value = 0.2773195847425811
$ws.Cells.Item(14, 17).Value = 2.622630400310666
$ws.Cells.Item(14, 18).Value = 23.603673602796
$ws.Cells.Item(14, 19).Value = 0.007521725750385039
$ws.Cells.Item(14, 20).Value = 0.007521725750385037

# Row 15
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.062692
$ws.Cells.Item(15, 8).Value = 0.188076
$ws.Cells.Item(15, 9).Value = 0.02712295187289783
$ws.Cells.Item(15, 10).Value = 0.02712295187289783
$ws.Cells.Item(15, 15).Value = 0.4239803668761465
$ws.Cells.Item(15, 16).Value = 0.4239803668761465
$ws.Cells.Item(15, 17).Value = 4.009611511341332
$ws.Cells.Item(15, 18).Value = 36.08650360207199
$ws.Cells.Item(15, 19).Value = 0.01149959908583529
$ws.Cells.Item(15, 20).Value = 0.01149959908583529

# Row 16
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.062692
$ws.Cells.Item(16, 8).Value = 0.188076
$ws.Cells.Item(16, 9).Value = 0.02712295187289783
$ws.Cells.Item(16, 10).Value = 0.02712295187289783
$ws.Cells.Item(16, 13).Value = 14.18032733333333
$ws.Cells.Item(16, 14).Value = 42.540982
$ws.Cells.Item(16, 15).Value = 0.09400302539123752
$ws.Cells.Item(16, 16).Value = 0.09400302539123752
$ws.Cells.Item(16, 17).Value = 0.8889930811813332
$ws.Cells.Item(16, 18).Value = 8.000937730632
$ws.Cells.Item(16, 19).Value = 0.002549639533593328
$ws.Cells.Item(16, 20).Value = 0.002549639533593328

# Row 17
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.062692
$ws.Cells.Item(17, 8).Value = 0.188076
$ws.Cells.Item(17, 9).Value = 0.02712295187289783
$ws.Cells.Item(17, 10).Value = 0.02712295187289783
$ws.Cells.Item(17, 13).Value = 30.87848266666667
$ws.Cells.Item(17, 14).Value = 92.635448
$ws.Cells.Item(17, 15).Value = 0.204697022990035
$ws.Cells.Item(17, 16).Value = 0.2046970229900349
$ws.Cells.Item(17, 17).Value = 1.935833835338667
$ws.Cells.Item(17, 18).Value = 17.422504518048
$ws.Cells.Item(17, 19).Value = 0.005551987503084179
$ws.Cells.Item(17, 20).Value = 0.005551987503084179

Write-Host "Updated cells for Bmp8a-Tgfbr2 TPM data"